$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the coordinate values Q2 (Ost) and R2 (Nord) to integers
$ws.Range("Q2").Value = 490412
$ws.Range("R2").Value = 6629327

# Remove the Starttid (Z2) and Sluttid (AB2) cell contents entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
